# Locate the target slide/shape by distinctive content rather than a
# hard-coded index, so the script is robust to any reordering.
$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null
$targetGroup = $null

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    for ($j = 1; $j -le $sl.Shapes.Count; $j++) {
        $shp = $sl.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            $shpText = $shp.TextFrame.TextRange.Text
            if ($shpText -like "*Consider the grammar for a while statement*") {
                $targetSlide = $sl
                $targetShape = $shp
            }
        }
        if ($shp.Name -eq "Group 1") {
            $targetGroup = $shp
        }
    }
}

$tf = $targetShape.TextFrame
$tr = $tf.TextRange

# ---------------------------------------------------------------------
# 1) "Example 2: Consider the grammar for a while statement." becomes
#    "Example 2: Consider the following rule for a " + "while" (Consolas)
#    + " statement."
# ---------------------------------------------------------------------
$introRun = $tr.Find("Example 2: Consider the grammar for a while statement.")
$introRun.Text = "Example 2: Consider the following rule for a while statement."

$full = $tf.TextRange
$whileInIntro = $full.Find("while", $introRun.Start)
$whileInIntro.Font.Name = "Consolas"

# ---------------------------------------------------------------------
# 2) ' = "while" "(" ' becomes ' = "while" '
# ---------------------------------------------------------------------
$full = $tf.TextRange
$ruleMiddle = $full.Find(' = "while" "(" ')
$ruleMiddle.Text = ' = "while" '

# ---------------------------------------------------------------------
# 3) ' ")" statement .' becomes ' "loop" statement .'
# ---------------------------------------------------------------------
$full = $tf.TextRange
$ruleTail = $full.Find(' ")" statement .')
$ruleTail.Text = ' "loop" statement .'

# ---------------------------------------------------------------------
# 4) Second paragraph: insert `"while"` and `"loop"` (Consolas) into the
#    sentence about terminal symbols.
# ---------------------------------------------------------------------
$full = $tf.TextRange
$p2Start = $full.Find("Once a while statement has been parsed")
$p2End = $full.Find("would contain only ")
$p2Len = ($p2End.Start + $p2End.Length) - $p2Start.Start
$p2Run = $full.Characters($p2Start.Start, $p2Len)

$apostrophe = [char]8217
$newP2Text = "Once a while statement has been parsed, we don" + $apostrophe + "t need to retain the terminal symbols " + [char]34 + "while" + [char]34 + " and " + [char]34 + "loop" + [char]34 + ".  The abstract syntax tree for a while statement would contain only "
$p2Run.Text = $newP2Text

$full = $tf.TextRange
$whileQuoted = $full.Find([char]34 + "while" + [char]34, $p2Start.Start)
$whileQuoted.Font.Name = "Consolas"

$full = $tf.TextRange
$loopQuoted = $full.Find([char]34 + "loop" + [char]34, $p2Start.Start)
$loopQuoted.Font.Name = "Consolas"

# ---------------------------------------------------------------------
# 5) Move the "Group 1" group shape down slightly (its Top offset
#    changes from 4191000 EMU to 4399848 EMU; 914400 EMU per inch,
#    12700 EMU per point).
# ---------------------------------------------------------------------
$targetGroup.Top = 4399848 / 12700
